$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5235
$ws.Range("J3").Value = 5558
$ws.Range("G4").Value = 1472
$ws.Range("H4").Value = 1698
$ws.Range("I4").Value = 1770
$ws.Range("J4").Value = 1234
$ws.Range("J5").Value = 434
$ws.Range("J6").Value = 6938
$ws.Range("G7").Value = 24696
$ws.Range("H7").Value = 26009
$ws.Range("I7").Value = 26222
$ws.Range("J7").Value = 19399

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 154
$ws.Range("J5").Value = 62
$ws.Range("J6").Value = 144
$ws.Range("J8").Value = 1229
$ws.Range("J9").Value = 96
$ws.Range("J10").Value = 132
$ws.Range("J11").Value = 304
$ws.Range("J14").Value = 97
$ws.Range("J15").Value = 212
$ws.Range("J18").Value = 164
$ws.Range("J19").Value = 555
$ws.Range("J20").Value = 403
$ws.Range("J22").Value = 52
$ws.Range("J23").Value = 187
$ws.Range("J24").Value = 60
$ws.Range("J26").Value = 44
$ws.Range("J27").Value = 115
$ws.Range("J29").Value = 1089
$ws.Range("H33").Value = 1303
$ws.Range("J33").Value = 885
$ws.Range("J36").Value = 269
$ws.Range("J37").Value = 608
$ws.Range("J41").Value = 123
$ws.Range("J42").Value = 793
$ws.Range("J48").Value = 225
$ws.Range("J52").Value = 489
$ws.Range("J53").Value = 257
$ws.Range("J54").Value = 370
$ws.Range("J55").Value = 254
$ws.Range("G63").Value = 270
$ws.Range("I63").Value = 235
$ws.Range("J63").Value = 66
$ws.Range("J64").Value = 129
$ws.Range("J66").Value = 61
$ws.Range("J67").Value = 749
$ws.Range("J68").Value = 38
$ws.Range("J75").Value = 58
$ws.Range("J78").Value = 244
$ws.Range("J79").Value = 556
$ws.Range("J84").Value = 161
$ws.Range("J85").Value = 829
$ws.Range("J89").Value = 248
$ws.Range("J94").Value = 189
$ws.Range("J96").Value = 234
$ws.Range("G101").Value = 24696
$ws.Range("H101").Value = 26009
$ws.Range("I101").Value = 26222
$ws.Range("J101").Value = 19399

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J3").Value = 62
$ws.Range("J4").Value = 14
$ws.Range("J6").Value = 84
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 119
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 79
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 216
$ws.Range("J3").Value = 305
$ws.Range("J6").Value = 236
$ws.Range("J7").Value = 829

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 145
$ws.Range("J6").Value = 202
$ws.Range("J7").Value = 489

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 341
$ws.Range("J3").Value = 374
$ws.Range("J5").Value = 35
$ws.Range("J6").Value = 404
$ws.Range("J7").Value = 1229

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 220
$ws.Range("J3").Value = 290
$ws.Range("H4").Value = 63
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 298
$ws.Range("H7").Value = 1303
$ws.Range("J7").Value = 885

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J5").Value = 25
$ws.Range("J6").Value = 174
$ws.Range("J7").Value = 608

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 185
$ws.Range("J4").Value = 59
$ws.Range("J7").Value = 749

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 51
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 72
$ws.Range("J6").Value = 174
$ws.Range("J7").Value = 370

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 330
$ws.Range("J3").Value = 371
$ws.Range("J6").Value = 286
$ws.Range("J7").Value = 1089

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 40
$ws.Range("J5").Value = 2
$ws.Range("J6").Value = 115
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 140
$ws.Range("J3").Value = 161
$ws.Range("J6").Value = 204
$ws.Range("J7").Value = 555

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 178
$ws.Range("J3").Value = 159
$ws.Range("J6").Value = 404
$ws.Range("J7").Value = 793

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J2").Value = 61
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 157
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 556

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 40
$ws.Range("J7").Value = 403

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 42
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 86
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 37
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J3").Value = 33
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J2").Value = 27
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 52
